$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 284, shifting rows 284:398 down to 285:399.
$ws.Rows("284:284").Insert()

# Populate the newly inserted row 284 with the new weekly data point.
# Columns that stay constant across this product's subset are copied
# from the row below (same values as every other row in this table).
$ws.Cells.Item(284, 1).Value = 3
$ws.Cells.Item(284, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(284, 3).Value = "Coquimbo"
$ws.Cells.Item(284, 4).Value = 44784
$ws.Cells.Item(284, 5).Value = 5
$ws.Cells.Item(284, 6).Value = 100112043
$ws.Cells.Item(284, 7).Value = "Pepino ensalada"
$ws.Cells.Item(284, 8).Value = "Sin especificar"
$ws.Cells.Item(284, 9).Value = "Primera"
$ws.Cells.Item(284, 10).Value = 85
$ws.Cells.Item(284, 11).Value = 17000
$ws.Cells.Item(284, 12).Value = 18000
$ws.Cells.Item(284, 13).Value = 17471
$ws.Cells.Item(284, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(284, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(284, 16).Value = 250
$ws.Cells.Item(284, 17).Value = 70
$ws.Cells.Item(284, 18).Value = "Hortaliza"
